$wb = $excel.ActiveWorkbook

# Rename the original first sheet from "sheet1" to "Sheet"
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Sheet"

# Clear out existing content on the first sheet (it becomes empty)
$ws1.Cells.Clear()

# Add three new sheets after the first one, each inserted right after the previous
$wsIntMapping = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$wsIntMapping.Name = "int_mapping"

$wsSystemName = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsIntMapping)
$wsSystemName.Name = "system_name"

$wsAddSuffix = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsSystemName)
$wsAddSuffix.Name = "add_suffix_system_name"

# ---- Populate int_mapping ----
$wsIntMapping.Range("A1").Value = "MVP1"
$wsIntMapping.Range("A2").Value = "['INT_L_EPP_EDW_RPTDDET','INT_L_EPP_EDW_RPTDDET_T1','INT_L_IM_FC_IMBCRM_1_0','INT_L_IM_IMACTM_A_0']"
$wsIntMapping.Range("A4").Value = "MVP2"
$wsIntMapping.Range("A5").Value = "['INT_E_LQ_GL_ACCT_SEGMENT_D','INT_E_LQ_GL_BAL_SUM_OC_D','INT_E_LQ_GL_BAL_SUM_RC_D','INT_E_LQ_GL_OC_SEGMENT_D','INT_E_LQ_GL_PROD_SEGMENT_D','INT_E_LQ_GL_RC_SEGMENT_D','INT_E_LQ_IM_ST_DAILY_D','INT_E_LQ_LN_BR_FOREIGN_MKT_D','INT_L_CCP_EDW_TXNALSB2K','INT_T_CCP_TXNALSB2K_01','INT_L_RCR_EVENTCODE','INT_T_RCR_EVENTCODE_01','INT_L_WLS_EDW_WLSPARAM','INT_L_WLS_EDW_WLSPROD','INT_T_WLS_WLSPROD_01']"
$wsIntMapping.Range("A7").Value = "MVP3"
$wsIntMapping.Range("A8").Value = "['INT_L_RBF_EDW_SDL_CLEANSING']"
$wsIntMapping.Range("A10").Value = "MVP4"
$wsIntMapping.Range("A11").Value = "['INT_L_BI_EDW_PARAM_ATS_COMP_CD']"
$wsIntMapping.Range("A13").Value = "MVP6"
$wsIntMapping.Range("A14").Value = "[]"

# ---- Populate system_name ----
$wsSystemName.Range("A1").Value = "MVP1"
$wsSystemName.Range("A2").Value = "ATM,GN,INV,LCS,MMS,MUREX,PLPS,PMH"
$wsSystemName.Range("A4").Value = "MVP2"
$wsSystemName.Range("A5").Value = "ATM,BCM,CCP,CFO,CHYO,EFS,IPS,LQ,PDPA,PMH,RCR,TMS,WLS"
$wsSystemName.Range("A7").Value = "MVP3"
$wsSystemName.Range("A8").Value = "ACTM,AM,BCC,CB,CCB,CIM,CSENT,CVA,ESN,FES,LEAD_UL,MRP,OLS,PRM,RBF,TRD,WSS"
$wsSystemName.Range("A10").Value = "MVP4"
$wsSystemName.Range("A11").Value = "BI,BIFI,EDW,ESL,SBG"
$wsSystemName.Range("A13").Value = "MVP6"
$wsSystemName.Range("A14").Value = "AML,CDD,CSM,ENL,LCS,OBM,RDS,SCBL"

# ---- Populate add_suffix_system_name ----
$wsAddSuffix.Range("A1").Value = "MVP1"
$wsAddSuffix.Range("A2").Value = "['REGISTER_CONFIG_SYSTEM_ATM','REGISTER_CONFIG_SYSTEM_GN','REGISTER_CONFIG_SYSTEM_INV','REGISTER_CONFIG_SYSTEM_LCS','REGISTER_CONFIG_SYSTEM_MMS','REGISTER_CONFIG_SYSTEM_MUREX','REGISTER_CONFIG_SYSTEM_PLPS','REGISTER_CONFIG_SYSTEM_PMH']"
$wsAddSuffix.Range("A4").Value = "MVP2"
$wsAddSuffix.Range("A5").Value = "['REGISTER_CONFIG_SYSTEM_ATM','REGISTER_CONFIG_SYSTEM_BCM','REGISTER_CONFIG_SYSTEM_CCP','REGISTER_CONFIG_SYSTEM_CFO','REGISTER_CONFIG_SYSTEM_CHYO','REGISTER_CONFIG_SYSTEM_EFS','REGISTER_CONFIG_SYSTEM_IPS','REGISTER_CONFIG_SYSTEM_LQ','REGISTER_CONFIG_SYSTEM_PDPA','REGISTER_CONFIG_SYSTEM_PMH','REGISTER_CONFIG_SYSTEM_RCR','REGISTER_CONFIG_SYSTEM_TMS','REGISTER_CONFIG_SYSTEM_WLS']"
$wsAddSuffix.Range("A7").Value = "MVP3"
$wsAddSuffix.Range("A8").Value = "['REGISTER_CONFIG_SYSTEM_ACTM','REGISTER_CONFIG_SYSTEM_AM','REGISTER_CONFIG_SYSTEM_BCC','REGISTER_CONFIG_SYSTEM_CB','REGISTER_CONFIG_SYSTEM_CCB','REGISTER_CONFIG_SYSTEM_CIM','REGISTER_CONFIG_SYSTEM_CSENT','REGISTER_CONFIG_SYSTEM_CVA','REGISTER_CONFIG_SYSTEM_ESN','REGISTER_CONFIG_SYSTEM_FES','REGISTER_CONFIG_SYSTEM_LEAD_UL','REGISTER_CONFIG_SYSTEM_MRP','REGISTER_CONFIG_SYSTEM_OLS','REGISTER_CONFIG_SYSTEM_PRM','REGISTER_CONFIG_SYSTEM_RBF','REGISTER_CONFIG_SYSTEM_TRD','REGISTER_CONFIG_SYSTEM_WSS']"
$wsAddSuffix.Range("A10").Value = "MVP4"
$wsAddSuffix.Range("A11").Value = "['REGISTER_CONFIG_SYSTEM_BI','REGISTER_CONFIG_SYSTEM_BIFI','REGISTER_CONFIG_SYSTEM_EDW','REGISTER_CONFIG_SYSTEM_ESL','REGISTER_CONFIG_SYSTEM_SBG']"
$wsAddSuffix.Range("A13").Value = "MVP6"
$wsAddSuffix.Range("A14").Value = "['REGISTER_CONFIG_SYSTEM_AML','REGISTER_CONFIG_SYSTEM_CDD','REGISTER_CONFIG_SYSTEM_CSM','REGISTER_CONFIG_SYSTEM_ENL','REGISTER_CONFIG_SYSTEM_LCS','REGISTER_CONFIG_SYSTEM_OBM','REGISTER_CONFIG_SYSTEM_RDS','REGISTER_CONFIG_SYSTEM_SCBL']"

$ws1.Activate()
